$wb = $excel.ActiveWorkbook

# --- Sheet "Test Cases" ---
$wsCases = $wb.Worksheets.Item("Test Cases")
$wsCases.Activate()
# Runmode for Login_01 changes from "no" to "Yes" (fixed input, so it can run again)
$wsCases.Range("C2").Value = "Yes"
$wsCases.Range("E7").Select()

# --- Sheet "Test Steps" ---
$wsSteps = $wb.Worksheets.Item("Test Steps")
$wsSteps.Activate()
# TS_019 result flips from PASS to FAIL
$wsSteps.Range("H20").Value = "FAIL"
# TS_020 .. TS_024 results get cleared out so the steps re-run cleanly
$wsSteps.Range("H21").ClearContents()
$wsSteps.Range("H22").ClearContents()
$wsSteps.Range("H23").ClearContents()
$wsSteps.Range("H24").ClearContents()
$wsSteps.Range("H25").ClearContents()
$wsSteps.Range("G8").Select()

# Test Steps remains the active/selected tab, matching the saved workbook state
$wsSteps.Activate()
